$d = $word.ActiveDocument

# Append two new bullet-list paragraphs after the current last paragraph
# ("IntelliJ is the most efficient IDE ..."). InsertParagraphAfter() on a
# collapsed end-of-document range clones the preceding paragraph's
# pPr/rPr (numPr, spacing, indent, justification, fonts, size, etc.), so
# the new paragraphs inherit the same bullet-list formatting.

$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last.Range
$p1.Collapse(0)
$p1.InsertAfter("SQL Server is a very competent implementation of SQL. SQL Server has a lot of functionality that other versions of SQL do not have. Furthermore, especially with stored procedures, SQL Server has more competent features in my opinion allowing for a faster and more secure implementation of our database than what would be possible in something like MySQL.")

$r2 = $d.Paragraphs.Last.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last.Range
$p2.Collapse(0)
$p2.InsertAfter("It’s better to agree on the format of and standardize the responses between the front end and backend.")
